# Auto-generated Excel COM-interop script
# Applies targeted cell updates (currentAveragePrice / LevePrice / LeveProfit columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# ALC row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 7408600.5
$ws.Range("I88").Value = 888.1111
$ws.Range("J88").Value = 18520168
$ws.Range("K88").Value = 888.1111
$ws.Range("L88").Value = 18520168
$ws.Range("M88").Value = -482.1111
$ws.Range("N88").Value = -18520980

# ALC row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 7408600.5
$ws.Range("I91").Value = 888.1111
$ws.Range("J91").Value = 18520168
$ws.Range("K91").Value = 888.1111
$ws.Range("L91").Value = 18520168
$ws.Range("M91").Value = 515.8889
$ws.Range("N91").Value = -18522976

# ALC row 95: Official Strategy Guide / Gyuki Leather Codex
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 500624
$ws.Range("J95").Value = 500624
$ws.Range("L95").Value = 500624
$ws.Range("N95").Value = -506116

# ALC row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 254856.55
$ws.Range("I132").Value = 296804.28
$ws.Range("J132").Value = 39874.5
$ws.Range("K132").Value = 890412.8400000001
$ws.Range("L132").Value = 119623.5
$ws.Range("M132").Value = -887882.8400000001
$ws.Range("N132").Value = -124683.5

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7694934
$ws.Range("I138").Value = 3242.2942
$ws.Range("J138").Value = 10419075
$ws.Range("K138").Value = 9726.882599999999
$ws.Range("L138").Value = 31257225
$ws.Range("M138").Value = -4586.882599999999
$ws.Range("N138").Value = -31267505

# ARM row 56: Feasting the Night Away / Hells' Kitchen
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 10552.941
$ws.Range("J56").Value = 10552.941
$ws.Range("L56").Value = 10552.941
$ws.Range("N56").Value = -12036.941

# ARM row 97: Ore for Me / High Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 22222726
$ws.Range("I97").Value = 25641514
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 25641514
$ws.Range("L97").Value = 600
$ws.Range("M97").Value = -25641018
$ws.Range("N97").Value = -1592

# ARM row 122: Haste for High Durium / High Durium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1036.9412
$ws.Range("I122").Value = 962.3333
$ws.Range("J122").Value = 1324.7142
$ws.Range("K122").Value = 2886.9999
$ws.Range("L122").Value = 3974.1426
$ws.Range("M122").Value = -436.9998999999998
$ws.Range("N122").Value = -8874.142599999999

# ARM row 128: Heading toward Bankruptcy / Manganese Helm of the Falling Dragon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 56851.668
$ws.Range("J128").Value = 56851.668
$ws.Range("L128").Value = 56851.668
$ws.Range("N128").Value = -66811.66800000001

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2625.025
$ws.Range("I132").Value = 1938.1072
$ws.Range("J132").Value = 4227.8335
$ws.Range("K132").Value = 5814.321599999999
$ws.Range("L132").Value = 12683.5005
$ws.Range("M132").Value = -3284.321599999999
$ws.Range("N132").Value = -17743.5005

# ARM row 133: Shielding My Students / Mountain Chromite Tower Shield
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 34125
$ws.Range("J133").Value = 34125
$ws.Range("L133").Value = 34125
$ws.Range("N133").Value = -39185

# ARM row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 94793
$ws.Range("J139").Value = 94793
$ws.Range("L139").Value = 94793
$ws.Range("N139").Value = -105073

# BSM row 59: Pop That Top / Cobalt Raising Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -51694

# BSM row 133: Paring Is Caring / Mountain Chromite Hatchet
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 38140
$ws.Range("I133").Value = 20000
$ws.Range("J133").Value = 44186.668
$ws.Range("K133").Value = 20000
$ws.Range("L133").Value = 44186.668
$ws.Range("M133").Value = -14940
$ws.Range("N133").Value = -54306.668

# BSM row 139: Maul Me / Titanium Gold Maul
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# CRP row 127: In Rod We Trust / Red Pine Fishing Rod
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 30865
$ws.Range("J127").Value = 30865
$ws.Range("L127").Value = 30865
$ws.Range("N127").Value = -40785

# CUL row 68: Such a Butter Face / Fermented Butter
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2083.9072
$ws.Range("I68").Value = 2480.106
$ws.Range("J68").Value = 1240.3871
$ws.Range("K68").Value = 7440.318000000001
$ws.Range("L68").Value = 3721.1613
$ws.Range("M68").Value = -6629.318000000001
$ws.Range("N68").Value = -5343.1613

# CUL row 71: No Margarine of Error (L) / Fermented Butter
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2083.9072
$ws.Range("I71").Value = 2480.106
$ws.Range("J71").Value = 1240.3871
$ws.Range("K71").Value = 22320.954
$ws.Range("L71").Value = 11163.4839
$ws.Range("M71").Value = -18264.954
$ws.Range("N71").Value = -19275.4839

# CUL row 129: Comfort Food / Yakow Moussaka
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1492.75
$ws.Range("J129").Value = 1584.8182
$ws.Range("L129").Value = 4754.4546
$ws.Range("N129").Value = -14754.4546

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3136.1865
$ws.Range("I131").Value = 501.7143
$ws.Range("J131").Value = 3490.827
$ws.Range("K131").Value = 1505.1429
$ws.Range("L131").Value = 10472.481
$ws.Range("M131").Value = 3534.8571
$ws.Range("N131").Value = -20552.481

# CUL row 137: Creative Chocolate / Gateau au Chocolat
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3370520
$ws.Range("I137").Value = 5885916
$ws.Range("J137").Value = 81156.16
$ws.Range("K137").Value = 17657748
$ws.Range("L137").Value = 243468.48
$ws.Range("M137").Value = -17652648
$ws.Range("N137").Value = -253668.48

# CUL row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 6169.696
$ws.Range("I140").Value = 7293.8667
$ws.Range("J140").Value = 4061.875
$ws.Range("K140").Value = 21881.6001
$ws.Range("L140").Value = 12185.625
$ws.Range("M140").Value = -16701.6001
$ws.Range("N140").Value = -22545.625

# GSM row 113: Copious Crystal Cannons / Manasilver Nugget
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1375.65
$ws.Range("I113").Value = 1197.5
$ws.Range("J113").Value = 1494.4166
$ws.Range("K113").Value = 1197.5
$ws.Range("L113").Value = 1494.4166
$ws.Range("M113").Value = 972.5
$ws.Range("N113").Value = -5834.4166

# GSM row 122: Awarding Academic Excellence / Ametrine
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1011927.6
$ws.Range("I122").Value = 1588429.2
$ws.Range("K122").Value = 4765287.6
$ws.Range("M122").Value = -4762837.6

# GSM row 123: Workplace Workout / Ametrine Ring of Fending
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10286.263
$ws.Range("J123").Value = 10286.263
$ws.Range("L123").Value = 10286.263
$ws.Range("N123").Value = -15186.263

# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2028.5555
$ws.Range("I126").Value = 1166.6666
$ws.Range("J126").Value = 2459.5
$ws.Range("K126").Value = 3499.9998
$ws.Range("L126").Value = 7378.5
$ws.Range("M126").Value = -1029.9998
$ws.Range("N126").Value = -12318.5

# GSM row 139: Ringing Gratitude / White Gold Ring of Healing
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 49999.5
$ws.Range("J139").Value = 49999.5
$ws.Range("L139").Value = 49999.5
$ws.Range("N139").Value = -60279.5

# LTW row 7: Tan Before the Ban / Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3500
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888

# LTW row 40: Best Served Toad / Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3339.111
$ws.Range("I40").Value = 3017.3333
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 3017.3333
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -2881.3333
$ws.Range("N40").Value = -3772

# LTW row 50: The Birdmen of Ishgard / Boarskin Culottes
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16274

# LTW row 61: Spelling Me Softly / Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2252.1304
$ws.Range("I61").Value = 2200.875
$ws.Range("J61").Value = 2369.2856
$ws.Range("K61").Value = 2200.875
$ws.Range("L61").Value = 2369.2856
$ws.Range("M61").Value = -1998.875
$ws.Range("N61").Value = -2773.2856

# LTW row 82: Trainin' the Neck / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 39884
$ws.Range("I82").Value = 101159.8
$ws.Range("J82").Value = 1586.625
$ws.Range("K82").Value = 101159.8
$ws.Range("L82").Value = 1586.625
$ws.Range("M82").Value = -100798.8
$ws.Range("N82").Value = -2308.625

# LTW row 85: Training Is Only Skintight (L) / Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 39884
$ws.Range("I85").Value = 101159.8
$ws.Range("J85").Value = 1586.625
$ws.Range("K85").Value = 101159.8
$ws.Range("L85").Value = 1586.625
$ws.Range("M85").Value = -99911.8
$ws.Range("N85").Value = -4082.625

# LTW row 113: Peace in Rest / Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2252.1304
$ws.Range("I113").Value = 2200.875
$ws.Range("J113").Value = 2369.2856
$ws.Range("K113").Value = 2200.875
$ws.Range("L113").Value = 2369.2856
$ws.Range("M113").Value = -30.875
$ws.Range("N113").Value = -6709.2856

# LTW row 126: Battered Books / Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

# WVR row 113: A Tender Table / Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1875.4166
$ws.Range("I113").Value = 2412.75
$ws.Range("K113").Value = 7238.25
$ws.Range("M113").Value = -5068.25

# WVR row 136: Weaving the Envelope / Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20897384
$ws.Range("I136").Value = 30394598
$ws.Range("J136").Value = 3513
$ws.Range("K136").Value = 91183794
$ws.Range("L136").Value = 10539
$ws.Range("M136").Value = -91181244
$ws.Range("N136").Value = -15639
